$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.482.74"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "1.693.79"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'316.45"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "'0.3911"
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("D8").Value = "'0.4045"
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").Value = "'1.479"
$ws.Range("E9").Value = "  -2.61%  "
$ws.Range("E10").Value = "  +0.03%  "
$ws.Range("D11").Value = "'53.34"
$ws.Range("E11").Value = "  +1.02%  "
$ws.Range("D12").Value = "'0.08778"
$ws.Range("E12").Value = "  -0.68%  "
$ws.Range("D13").Value = "'26.12"
$ws.Range("E13").Value = "  +10.84%  "
$ws.Range("D14").Value = "'7.429"
$ws.Range("E14").Value = "  +0.35%  "
$ws.Range("D15").Value = "'8.126"
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("D16").Value = "'0.00001352"
$ws.Range("E16").Value = "  +2.47%  "
$ws.Range("D17").Value = "1.694.34"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("D18").Value = "'97.54"
$ws.Range("E18").Value = "  -1.93%  "
$ws.Range("D19").Value = "'0.07164"
$ws.Range("E19").Value = "  +2.12%  "
$ws.Range("D20").Value = "'20.39"
$ws.Range("E20").Value = "  +3.23%  "
$ws.Range("D21").Value = "'7.286"
$ws.Range("E21").Value = "  +3.10%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  -0.37%  "
$ws.Range("D23").Value = "'14.26"
$ws.Range("E23").Value = "  -2.44%  "
$ws.Range("D24").Value = "24.492.20"
$ws.Range("E24").Value = "  -0.63%  "
$ws.Range("D25").Value = "'2.988"
$ws.Range("E25").Value = "  -4.68%  "
$ws.Range("D26").Value = "'2.329"
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("D27").Value = "'22.88"
$ws.Range("E27").Value = "  +1.14%  "
$ws.Range("D28").Value = "'169.03"
$ws.Range("E28").Value = "  +3.64%  "
$ws.Range("D29").Value = "'5.816"
$ws.Range("E29").Value = "  +12.99%  "
$ws.Range("D30").Value = "'144.55"
$ws.Range("E30").Value = "  +6.79%  "
$ws.Range("D31").Value = "'8.366"
$ws.Range("E31").Value = "  -4.02%  "
$ws.Range("D32").Value = "1.880.95"
$ws.Range("E32").Value = "  -0.41%  "
$ws.Range("D33").Value = "'2.183"
$ws.Range("E33").Value = "  +11.14%  "
$ws.Range("D34").Value = "'0.08751"
$ws.Range("E34").Value = "  -2.16%  "
$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D35").Value = "'0.03153"
$ws.Range("E35").Value = "  +9.06%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'1.040"
$ws.Range("E36").Value = "  -2.45%  "
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = "'7.174"
$ws.Range("E37").Value = "  -5.43%  "
$ws.Range("D38").Value = "'0.2799"
$ws.Range("E38").Value = "  +1.92%  "
$ws.Range("D39").Value = "'10.85"
$ws.Range("E39").Value = "  -1.86%  "
$ws.Range("D40").Value = "'0.09156"
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("D41").Value = "'14.12"
$ws.Range("E41").Value = "  -1.97%  "
$ws.Range("D42").Value = "'0.7995"
$ws.Range("E42").Value = "  +4.60%  "
$ws.Range("D43").Value = "'1.480"
$ws.Range("E43").Value = "  +1.65%  "
$ws.Range("D44").Value = "'17.26"
$ws.Range("E44").Value = "  +9.30%  "
$ws.Range("D45").Value = "'2.647"
$ws.Range("E45").Value = "  +3.52%  "
$ws.Range("D46").Value = "'0.7243"
$ws.Range("E46").Value = "  +1.14%  "
$ws.Range("D47").Value = "'4.251"
$ws.Range("E47").Value = "  +0.87%  "
$ws.Range("D48").Value = "'1.380"
$ws.Range("E48").Value = "  +3.35%  "
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").Value = "'139.64"
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("D51").Value = "'0.08201"
$ws.Range("E51").Value = "  +2.98%  "
